$wb = $excel.ActiveWorkbook

# --- Sheet "Info": update Objetivo / Tiempo result row ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("A2").Value = 640108574274.0112
$wsInfo.Range("B2").Value = 2.354000091552734

# --- Sheet "Activados": replace the 3-row sample with 19 rows (Proceso=1, Tiempo 0..360 step 20) ---
$wsActivados = $wb.Worksheets.Item("Activados")
for ($i = 0; $i -le 18; $i++) {
    $row = 2 + $i
    $wsActivados.Cells.Item($row, 1).Value = 1
    $wsActivados.Cells.Item($row, 2).Value = $i * 20
}

# --- Sheet "Operando": Proceso column changes from 4 to 1 for every data row (Tiempo stays the same) ---
$wsOperando = $wb.Worksheets.Item("Operando")
$wsOperando.Range("A2:A366").Value = 1

# --- Sheet "Contaminantes": refreshed B (Z) / C (Concentracion) values ---
$wsContaminantes = $wb.Worksheets.Item("Contaminantes")
$wsContaminantes.Range("B2").Value = 449208244800.0004
$wsContaminantes.Range("C2").Value = 16.66000000000001
$wsContaminantes.Range("B3").Value = 13481640000.00001
$wsContaminantes.Range("C3").Value = 0.5000000000000004
$wsContaminantes.Range("B4").Value = 87091394399.99998
$wsContaminantes.Range("C4").Value = 3.23
$wsContaminantes.Range("B5").Value = 307074.010608
$wsContaminantes.Range("C5").Value = 0.0000113886
$wsContaminantes.Range("B6").Value = 90326988000.00008
$wsContaminantes.Range("C6").Value = 3.350000000000003
